$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
$ws.Range("A7").Value = "a1"
$ws.Range("B7").Value = 0.8752631545066833
$ws.Range("C7").Value = 0.5897498726844788
$ws.Range("D7").Value = 0.837644100189209
$ws.Range("E7").Value = 0.8836870193481445
$ws.Range("F7").Value = 0.8707107901573181
$ws.Range("G7").Value = 148.128173828125
$ws.Range("H7").Value = 13.14312839508057
$ws.Range("I7").Value = 13.35710144042969
$ws.Range("J7").Value = 86.67254638671875
$ws.Range("K7").Value = 91.16632843017578

# --- Row 8 ---
$ws.Range("A8").Value = "b2"
$ws.Range("B8").Value = 0.875263512134552
$ws.Range("C8").Value = 0.5897493362426758
$ws.Range("D8").Value = 0.8376448750495911
$ws.Range("E8").Value = 0.8836870193481445
$ws.Range("F8").Value = 0.8707107305526733
$ws.Range("G8").Value = 148.1097259521484
$ws.Range("H8").Value = 13.14150238037109
$ws.Range("I8").Value = 13.35542774200439
$ws.Range("J8").Value = 86.67243957519531
$ws.Range("K8").Value = 91.16508483886719

# Match the formatting used by the other column-A label cells (bold,
# bordered, centered/top-aligned) by copying formats from A6.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
